# Auto commit at 2025-10-29  8:03:27.26
# Refresh the "Metrics" sheet's raw totals with the next day's figures and
# move the saved cell-selection on both the "Metrics" and "today" sheets.
# Downstream formulas on "today" (and the E/F helper columns) recompute
# automatically from these Metrics values when Excel recalculates.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 389705.19
$metrics.Range("B3").Value  = 319547.86000000004
$metrics.Range("B4").Value  = 124064.73999999999
$metrics.Range("B5").Value  = 15578
$metrics.Range("B6").Value  = 4756836.66
$metrics.Range("B7").Value  = 4009366.5299999993
$metrics.Range("B8").Value  = 1394666.88
$metrics.Range("B9").Value  = 184579
$metrics.Range("B10").Value = 33222160.460999824
$metrics.Range("B11").Value = 31284588.050000001
$metrics.Range("B12").Value = 11676375.770000005
$metrics.Range("B13").Value = 1282206

# Move the stored selection on "Metrics" from E20 to F21.
$metrics.Activate() | Out-Null
$metrics.Range("F21").Select() | Out-Null

# Move the stored selection on "today" from D8 to G8.
$today = $wb.Worksheets.Item("today")
$today.Activate() | Out-Null
$today.Range("G8").Select() | Out-Null
